$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..47 (columns A..G)
# A=Id (number), B=IP, C=Dispositivo, D=Puerto, E=status, F=Unnamed:4, G=Unnamed:5
$data = @(
    ,@(1277, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/0", "nan", "nan", "PUERTOLIBRE")
    ,@(1321, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/30", "nan", "nan", "PUERTOLIBRE")
    ,@(1322, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/31", "nan", "nan", "PUERTOLIBRE")
    ,@(1323, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/32", "nan", "nan", "PUERTOLIBRE")
    ,@(1324, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/33", "nan", "nan", "PUERTOLIBRE")
    ,@(1325, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/34", "nan", "nan", "PUERTOLIBRE")
    ,@(1326, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/35", "nan", "nan", "PUERTOLIBRE")
    ,@(1327, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/36", "nan", "nan", "PUERTOLIBRE")
    ,@(1328, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/37", "nan", "nan", "PUERTOLIBRE")
    ,@(1329, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/38", "nan", "nan", "PUERTOLIBRE")
    ,@(1330, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/39", "nan", "nan", "PUERTOLIBRE")
    ,@(1332, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/40", "nan", "nan", "PUERTOLIBRE")
    ,@(1333, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/41", "nan", "nan", "PUERTOLIBRE")
    ,@(1334, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/42", "nan", "nan", "PUERTOLIBRE")
    ,@(1335, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/43", "nan", "nan", "PUERTOLIBRE")
    ,@(1336, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/44", "nan", "nan", "PUERTOLIBRE")
    ,@(1337, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/45", "nan", "nan", "PUERTOLIBRE")
    ,@(1338, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/46", "nan", "nan", "PUERTOLIBRE")
    ,@(1339, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/47", "nan", "nan", "PUERTOLIBRE")
    ,@(1340, "172.24.254.23", "MEDE-CABA-H-01-DAAS", "xe-0/0/48", "nan", "nan", "PUERTOLIBRE")
    ,@(1347, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/0", "nan", "nan", "PUERTOLIBRE")
    ,@(1399, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/34", "nan", "nan", "PUERTOLIBRE")
    ,@(1414, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/41", "nan", "nan", "PUERTOLIBRE")
    ,@(1415, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/42", "nan", "nan", "PUERTOLIBRE")
    ,@(1416, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/43", "nan", "nan", "PUERTOLIBRE")
    ,@(1417, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/44", "nan", "nan", "PUERTOLIBRE")
    ,@(1418, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/45", "nan", "nan", "PUERTOLIBRE")
    ,@(1419, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/46", "nan", "nan", "PUERTOLIBRE")
    ,@(1420, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/47", "nan", "nan", "PUERTOLIBRE")
    ,@(1421, "172.24.254.24", "MEDE-CABA-H-02-DAAS", "xe-0/0/48", "nan", "nan", "PUERTOLIBRE")
    ,@(1431, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/0", "nan", "nan", "PUERTOLIBRE")
    ,@(1487, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/36", "nan", "nan", "PUERTOLIBRE")
    ,@(1489, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/37", "nan", "nan", "PUERTOLIBRE")
    ,@(1491, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/38", "nan", "nan", "PUERTOLIBRE")
    ,@(1493, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/39", "nan", "nan", "PUERTOLIBRE")
    ,@(1496, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/40", "nan", "nan", "PUERTOLIBRE")
    ,@(1498, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/41", "nan", "nan", "PUERTOLIBRE")
    ,@(1499, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/42", "nan", "nan", "PUERTOLIBRE")
    ,@(1500, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/43", "nan", "nan", "PUERTOLIBRE")
    ,@(1501, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/44", "nan", "nan", "PUERTOLIBRE")
    ,@(1502, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/45", "nan", "nan", "PUERTOLIBRE")
    ,@(1503, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/46", "nan", "nan", "PUERTOLIBRE")
    ,@(1504, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/47", "nan", "nan", "PUERTOLIBRE")
    ,@(1505, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/48", "nan", "nan", "PUERTOLIBRE")
    ,@(1507, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/5", "nan", "nan", "PUERTOLIBRE")
    ,@(1509, "172.24.254.25", "MEDE-CABA-H-03-DAAS", "xe-0/0/7", "nan", "nan", "PUERTOLIBRE")
)

$startRow = 2
$origLastRow = 35
$newLastRow = $startRow + $data.Count - 1

# For rows beyond the originally populated range, copy the formatting
# (style) from row 2's column A cell so new rows match the existing
# pattern (bordered, bold, centered numeric id cells).
if ($newLastRow -gt $origLastRow) {
    $ws.Cells.Item(2, 1).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($origLastRow + 1, 1), $ws.Cells.Item($newLastRow, 1)).PasteSpecial(-4122) | Out-Null
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
}
